$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 11:52"

# Update data values per-row (B=Casos totales, C=Casos activos, D=Recuperados, E=Muertes)
$updates = @{
    4  = @{ B = 60765; C = 36314; D = 16346; E = 8105 }
    5  = @{ B = 48654; C = 19060; D = 24689; E = 4905 }
    6  = @{ B = 16690; C = 6448;  D = 8506;  E = 1736 }
    7  = @{ B = 15785; C = 5382;  D = 7967;  E = 2436 }
    8  = @{ B = 12619; C = 10474; D = 871;   E = 1274 }
    9  = @{ B = 12004; C = 5272;  D = 5544;  E = 1188 }
    10 = @{ B = 9466;  C = 5393;  D = 3526;  E = 547 }
    14 = @{ B = 5042;  C = 2189;  D = 2117;  E = 736 }
    15 = @{ B = 4794;  C = 2082;  D = 2264;  E = 448 }
    17 = @{ B = 3910;  C = 2091;  D = 1489;  E = 330 }
    23 = @{ B = 2764;  C = 1797;  D = 527;   E = 440 }
    30 = @{ B = 2266;  C = 809;   D = 1191;  E = 266 }
    32 = @{ B = 2202;  C = 1107;  D = 961;   E = 134 }
    33 = @{ B = 2146;  C = 1372;  D = 583;   E = 191 }
    38 = @{ B = 1480;  C = 1153;  D = 197 }
    59 = @{ C = 95; D = 17 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
